# Fruta / hortaliza, semanal
#
# The weekly refresh reshuffles the data rows (2-20) of the single sheet:
# each row keeps its "identity" columns (A,B,C,E-K) but the
# date/quality/volume/price/unit/origin columns (D, L-T) are redistributed
# among the rows according to a fixed permutation. Capture the original
# values first, then write them back in the new order so overlapping
# writes do not clobber data we still need to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get reshuffled: D, L, M, N, O, P, Q, R, S, T
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot the current contents of the permuted columns for every data row.
$original = @{}
foreach ($r in 2..20) {
    $rowValues = @{}
    foreach ($c in $cols) {
        $rowValues[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowValues
}

# Target row -> source row (the row whose old D/L-T values now belong here).
$map = @{
    2  = 15
    3  = 2
    4  = 20
    5  = 13
    6  = 19
    7  = 7
    8  = 16
    9  = 18
    10 = 8
    11 = 14
    12 = 17
    13 = 5
    14 = 9
    15 = 12
    16 = 4
    17 = 10
    18 = 11
    19 = 6
    20 = 3
}

foreach ($target in 2..20) {
    $source = $map[$target]
    $src = $original[$source]
    foreach ($c in $cols) {
        $ws.Cells.Item($target, $c).Value = $src[$c]
    }
}
